# Aulas modulo 2 - adiciona Aula 2 e Aula 4 apos a Aula 1 existente.
$d = $word.ActiveDocument

# Paragraph 9 currently holds "Access Passcode: *v27V&As" - the new
# content must be inserted right after it (and before the trailing
# empty paragraph that closes the document).
$anchorIndex = 9
$anchor = $d.Paragraphs($anchorIndex)
Write-Host "Anchor paragraph text: " $anchor.Range.Text

$lines = @(
  "#########",
  "Aula 2:",
  "Topic: 780 - Data Science Degree",
  "Start Time : Aug 4, 2021 06:45 PM",
  "",
  "Meeting Recording:",
  "https://us02web.zoom.us/rec/share/PysSjuYgnrTJMWBrbL4_5Tp2raBK0szxn5sT8-0t0I4QJI8xSfN9RHlfB7q0_2Ax.Lt6qqX96k_ePasev",
  "",
  "Access Passcode: 7d3K!xit",
  "",
  "#######",
  "Aula 4:",
  "Topic: 780 - Data Science Degree",
  "Start Time : Aug 9, 2021 06:45 PM",
  "",
  "Meeting Recording:",
  "https://us02web.zoom.us/rec/share/rmozQk8bWImbNlyMcBbU-OsWrnFa_WfYaVMMt13h4tJ6Aefsqx6Ktn3OJuo_-JPB.2wEQD8HeSZit_EGT",
  "",
  "Access Passcode: 2d#0zc22"
)

$idx = $anchorIndex
foreach ($line in $lines) {
  $p = $d.Paragraphs($idx)
  $p.Range.InsertParagraphAfter()
  $idx = $idx + 1
  if ($line -ne "") {
    $newPara = $d.Paragraphs($idx)
    $newPara.Range.InsertBefore($line)
  }
}

Write-Host "Final paragraph count: " $d.Paragraphs.Count
